# Bot-Varredura/dominios.xlsx — add MUNICIPIO/ESTADO columns, refresh a
# couple of timestamps, and append two new scanned domains (rows 93-94).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells D1 (MUNICIPIO) / E1 (ESTADO) -------------------------
# Copy C1's format (bold/centered/bordered header style) onto D1/E1 so we
# reuse the existing header style instead of minting new ones, then set text.
$ws.Range("C1").Copy($ws.Range("D1"))
$ws.Range("C1").Copy($ws.Range("E1"))
$ws.Range("D1").Value = "MUNICIPIO"
$ws.Range("E1").Value = "ESTADO"

# --- Row 2: Nova Friburgo / RJ — refreshed extraction timestamp ------------
$ws.Range("C2").Value = "2025-01-21 17:03:53"
$ws.Range("D2").Value = "Nova Friburgo"
$ws.Range("E2").Value = "Rio de Janeiro"

# --- Row 3: Bom Jardim / RJ — refreshed extraction timestamp ---------------
$ws.Range("C3").Value = "2025-01-21 12:57:19"
$ws.Range("D3").Value = "Bom Jardim"
$ws.Range("E3").Value = "Rio de Janeiro"

# --- Rows 4-92: new MUNICIPIO/ESTADO columns exist but are still blank -----
# Touch each cell with a formatting no-op (LineStyle already "none") so the
# cell is materialised (extends dimension / row spans) without changing its
# appearance or allocating a new cell style.
$ws.Range("D4:E92").Borders.LineStyle = -4142

# --- New rows: domains scanned on 2025-01-21 --------------------------------
$ws.Range("A93").Value = "acrelandia.ac.gov.br"
$ws.Range("B93").Value = "SUCESSO"
$ws.Range("C93").Value = "2025-01-21 17:04:11"
$ws.Range("D93").Value = "Acrelândia"
$ws.Range("E93").Value = "Acre"

$ws.Range("A94").Value = "assisbrasil.ac.gov.br"
$ws.Range("B94").Value = "SUCESSO"
$ws.Range("C94").Value = "2025-01-21 12:57:32"
$ws.Range("D94").Value = "Assis Brasil"
$ws.Range("E94").Value = "Acre"
